$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.009860666666667
$ws.Range("H2").Value = 3.029582
$ws.Range("I2").Value = 0.01353413605720072
$ws.Range("J2").Value = 0.01542521070970148
$ws.Range("M2").Value = 1.135186
$ws.Range("N2").Value = 3.405558
$ws.Range("O2").Value = 0.006415563145489937
$ws.Range("P2").Value = 0.00646579730555003
$ws.Range("Q2").Value = 1.146379690750667
$ws.Range("R2").Value = 10.317417216756
$ws.Range("S2").Value = 0.00008682910449462344
$ws.Range("T2").Value = 0.0000997362858443293
$ws.Range("G3").Value = 1.009860666666667
$ws.Range("H3").Value = 3.029582
$ws.Range("I3").Value = 0.01353413605720072
$ws.Range("J3").Value = 0.01542521070970148
$ws.Range("O3").Value = 0.02672480471352731
$ws.Range("P3").Value = 0.02693406118674866
$ws.Range("Q3").Value = 4.775383340183112
$ws.Range("R3").Value = 42.978450061648
$ws.Range("S3").Value = 0.0003616971430949978
$ws.Range("T3").Value = 0.0004154635690735904
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 1.009860666666667
$ws.Range("H4").Value = 3.029582
$ws.Range("I4").Value = 0.01353413605720072
$ws.Range("J4").Value = 0.01542521070970148
$ws.Range("M4").Value = 84.55360633333333
$ws.Range("N4").Value = 253.660819
$ws.Range("O4").Value = 0.4778591355164685
$ws.Range("P4").Value = 0.4816007949398642
$ws.Range("Q4").Value = 85.38736126085089
$ws.Range("R4").Value = 768.486251347658
$ws.Range("S4").Value = 0.006467410556256203
$ws.Range("T4").Value = 0.007428793739907139
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 1.009860666666667
$ws.Range("H5").Value = 3.029582
$ws.Range("I5").Value = 0.01353413605720072
$ws.Range("J5").Value = 0.01542521070970148
$ws.Range("M5").Value = 4.124113
$ws.Range("N5").Value = 8.248226000000001
$ws.Range("O5").Value = 0.02330764066032874
$ws.Range("P5").Value = 0.01566009371925767
$ws.Range("Q5").Value = 4.164779503588667
$ws.Range("R5").Value = 24.988677021532
$ws.Range("S5").Value = 0.0003154487798692329
$ws.Range("T5").Value = 0.0002415602453532223
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 1.009860666666667
$ws.Range("H6").Value = 3.029582
$ws.Range("I6").Value = 0.01353413605720072
$ws.Range("J6").Value = 0.01542521070970148
$ws.Range("M6").Value = 82.400874
$ws.Range("N6").Value = 247.202622
$ws.Range("O6").Value = 0.4656928559641855
$ws.Range("P6").Value = 0.4693392528485795
$ws.Range("Q6").Value = 83.213401551556
$ws.Range("R6").Value = 748.920613964004
$ws.Range("S6").Value = 0.006302750473485666
$ws.Range("T6").Value = 0.0072396568695232
$ws.Range("I7").Value = 0.6185519418990597
$ws.Range("J7").Value = 0.704979911415303
$ws.Range("M7").Value = 1.135186
$ws.Range("N7").Value = 3.405558
$ws.Range("O7").Value = 0.006415563145489937
$ws.Range("P7").Value = 0.00646579730555003
$ws.Range("Q7").Value = 52.39310295615066
$ws.Range("R7").Value = 471.5379266053559
$ws.Range("S7").Value = 0.00396835904201884
$ws.Range("T7").Value = 0.004558257211695965
$ws.Range("I8").Value = 0.6185519418990597
$ws.Range("J8").Value = 0.704979911415303
$ws.Range("O8").Value = 0.02672480471352731
$ws.Range("P8").Value = 0.02693406118674866
$ws.Range("S8").Value = 0.01653067985242546
$ws.Range("T8").Value = 0.01898797206948842
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("I9").Value = 0.6185519418990597
$ws.Range("J9").Value = 0.704979911415303
$ws.Range("M9").Value = 84.55360633333333
$ws.Range("N9").Value = 253.660819
$ws.Range("O9").Value = 0.4778591355164685
$ws.Range("P9").Value = 0.4816007949398642
$ws.Range("Q9").Value = 3902.466910212217
$ws.Range("R9").Value = 35122.20219190996
$ws.Range("S9").Value = 0.2955806962279175
$ws.Range("T9").Value = 0.339518885754245
$ws.Range("D10").Value = "MuSCs"
$ws.Range("I10").Value = 0.6185519418990597
$ws.Range("J10").Value = 0.704979911415303
$ws.Range("M10").Value = 4.124113
$ws.Range("N10").Value = 8.248226000000001
$ws.Range("O10").Value = 0.02330764066032874
$ws.Range("P10").Value = 0.01566009371925767
$ws.Range("Q10").Value = 190.3433243642887
$ws.Range("R10").Value = 1142.059946185732
$ws.Range("S10").Value = 0.01441698639153182
$ws.Range("T10").Value = 0.01104005148295762
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("I11").Value = 0.6185519418990597
$ws.Range("J11").Value = 0.704979911415303
$ws.Range("M11").Value = 82.400874
$ws.Range("N11").Value = 247.202622
$ws.Range("O11").Value = 0.4656928559641855
$ws.Range("P11").Value = 0.4693392528485795
$ws.Range("Q11").Value = 3803.110217320156
$ws.Range("R11").Value = 34227.9919558814
$ws.Range("S11").Value = 0.2880552203851661
$ws.Range("T11").Value = 0.3308747448969161
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.009315666666666667
$ws.Range("H12").Value = 0.027947
$ws.Range("I12").Value = 0.0001248484115599408
$ws.Range("J12").Value = 0.000142293017222847
$ws.Range("M12").Value = 1.135186
$ws.Range("N12").Value = 3.405558
$ws.Range("O12").Value = 0.006415563145489937
$ws.Range("P12").Value = 0.00646579730555003
$ws.Range("Q12").Value = 0.01057501438066667
$ws.Range("R12").Value = 0.095175129426
$ws.Range("S12").Value = 0.000000800972867976916
$ws.Range("T12").Value = 0.0000009200378073580682
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.009315666666666667
$ws.Range("H13").Value = 0.027947
$ws.Range("I13").Value = 0.0001248484115599408
$ws.Range("J13").Value = 0.000142293017222847
$ws.Range("O13").Value = 0.02672480471352731
$ws.Range("P13").Value = 0.02693406118674866
$ws.Range("Q13").Value = 0.04405150222311111
$ws.Range("R13").Value = 0.3964635200080001
$ws.Range("S13").Value = 0.000003336549417733504
$ws.Range("T13").Value = 0.000003832528832327242
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.009315666666666667
$ws.Range("H14").Value = 0.027947
$ws.Range("I14").Value = 0.0001248484115599408
$ws.Range("J14").Value = 0.000142293017222847
$ws.Range("M14").Value = 84.55360633333333
$ws.Range("N14").Value = 253.660819
$ws.Range("O14").Value = 0.4778591355164685
$ws.Range("P14").Value = 0.4816007949398642
$ws.Range("Q14").Value = 0.7876732120658889
$ws.Range("R14").Value = 7.089058908593
$ws.Range("S14").Value = 0.00005965995401863759
$ws.Range("T14").Value = 0.0000685284302089149
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.009315666666666667
$ws.Range("H15").Value = 0.027947
$ws.Range("I15").Value = 0.0001248484115599408
$ws.Range("J15").Value = 0.000142293017222847
$ws.Range("M15").Value = 4.124113
$ws.Range("N15").Value = 8.248226000000001
$ws.Range("O15").Value = 0.02330764066032874
$ws.Range("P15").Value = 0.01566009371925767
$ws.Range("Q15").Value = 0.03841886200366667
$ws.Range("R15").Value = 0.230513172022
$ws.Range("S15").Value = 0.000002909921913651933
$ws.Range("T15").Value = 0.00000222832198530573
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.009315666666666667
$ws.Range("H16").Value = 0.027947
$ws.Range("I16").Value = 0.0001248484115599408
$ws.Range("J16").Value = 0.000142293017222847
$ws.Range("M16").Value = 82.400874
$ws.Range("N16").Value = 247.202622
$ws.Range("O16").Value = 0.4656928559641855
$ws.Range("P16").Value = 0.4693392528485795
$ws.Range("Q16").Value = 0.767619075226
$ws.Range("R16").Value = 6.908571677034001
$ws.Range("S16").Value = 0.00005814101334194087
$ws.Range("T16").Value = 0.00006678369838894107
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 27.4428835
$ws.Range("H17").Value = 54.885767
$ws.Range("I17").Value = 0.3677890736321797
$ws.Range("J17").Value = 0.2794525848577725
$ws.Range("M17").Value = 1.135186
$ws.Range("N17").Value = 3.405558
$ws.Range("O17").Value = 0.006415563145489937
$ws.Range("P17").Value = 0.00646579730555003
$ws.Range("Q17").Value = 31.152777148831
$ws.Range("R17").Value = 186.916662892986
$ws.Range("S17").Value = 0.002359574026108497
$ws.Range("T17").Value = 0.001806883770202377
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 27.4428835
$ws.Range("H18").Value = 54.885767
$ws.Range("I18").Value = 0.3677890736321797
$ws.Range("J18").Value = 0.2794525848577725
$ws.Range("O18").Value = 0.02672480471352731
$ws.Range("P18").Value = 0.02693406118674866
$ws.Range("Q18").Value = 129.7706634174147
$ws.Range("R18").Value = 778.6239805044881
$ws.Range("S18").Value = 0.00982909116858912
$ws.Range("T18").Value = 0.007526793019354318
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 27.4428835
$ws.Range("H19").Value = 54.885767
$ws.Range("I19").Value = 0.3677890736321797
$ws.Range("J19").Value = 0.2794525848577725
$ws.Range("M19").Value = 84.55360633333333
$ws.Range("N19").Value = 253.660819
$ws.Range("O19").Value = 0.4778591355164685
$ws.Range("P19").Value = 0.4816007949398642
$ws.Range("Q19").Value = 2320.394768110529
$ws.Range("R19").Value = 13922.36860866317
$ws.Range("S19").Value = 0.1757513687782762
$ws.Range("T19").Value = 0.1345845870155031
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 27.4428835
$ws.Range("H20").Value = 54.885767
$ws.Range("I20").Value = 0.3677890736321797
$ws.Range("J20").Value = 0.2794525848577725
$ws.Range("M20").Value = 4.124113
$ws.Range("N20").Value = 8.248226000000001
$ws.Range("O20").Value = 0.02330764066032874
$ws.Range("P20").Value = 0.01566009371925767
$ws.Range("Q20").Value = 113.1775525998355
$ws.Range("R20").Value = 452.710210399342
$ws.Range("S20").Value = 0.008572295567014032
$ws.Range("T20").Value = 0.004376253668961525
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 27.4428835
$ws.Range("H21").Value = 54.885767
$ws.Range("I21").Value = 0.3677890736321797
$ws.Range("J21").Value = 0.2794525848577725
$ws.Range("M21").Value = 82.400874
$ws.Range("N21").Value = 247.202622
$ws.Range("O21").Value = 0.4656928559641855
$ws.Range("P21").Value = 0.4693392528485795
$ws.Range("Q21").Value = 2261.317585480179
$ws.Range("R21").Value = 13567.90551288108
$ws.Range("S21").Value = 0.1712767440921919
$ws.Range("T21").Value = 0.1311580673837512
